$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the header row (row 1) entirely
$ws.Range("A1:H1").ClearContents()

# Update row 2 values
$ws.Range("B2").Value = "JJ "
$ws.Range("C2").Value = "4/23/2025"

# Update row 3 values
$ws.Range("A3").Value = "Matthew wolz"
$ws.Range("B3").Value = "JJ"
$ws.Range("C3").Value = "4/23/2025"
$ws.Range("E3").Value = 1

# Remove rows 4, 5, 6 entirely
$ws.Range("A4:H6").ClearContents()
